$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 15
$ws.Range("I2").Value = -1.492114695340502
$ws.Range("J2").Value = -0.5420694593626935
$ws.Range("K2").Value = 0.2988218493395215

# Row 5 updates
$ws.Range("D5").Value = 2
$ws.Range("I5").Value = 2.908960573476702
$ws.Range("J5").Value = -3.248120300751879
$ws.Range("K5").Value = 2.272045697965013

# Row 6 updates
$ws.Range("D6").Value = 2
$ws.Range("I6").Value = 0.9738351254480289
$ws.Range("J6").Value = -2.152882205513786
$ws.Range("K6").Value = 2.130310603355944
